# Actualizacion Datos Personales 4 nov
# Updates Aprobados/Reprobados/Por_Apro/Por_Repro/Promedio/Blancos/Por_Blan
# for rows 2-7 (5ALCM, 5APM, 5ARHM, 5BEM, 5BLCM, 5ARHV) on sheets
# "1er Parcial" and "3er Parcial". "2o Parcial" is left untouched.

$wb = $excel.ActiveWorkbook

# Row => E (Aprobados), F (Reprobados), G (Por_Apro), H (Por_Repro),
#        I (Promedio), J (Blancos), K (Por_Blan)
$updates = @{
    2 = @(30, 5, 85.70999999999999, 14.29, 8,   0, 0)
    3 = @(18, 8, 69.23,             30.77, 6.2, 1, 3.85)
    4 = @(23, 0, 100,               0,     8,   0, 0)
    5 = @(20, 12, 62.5,             37.5,  6.9, 4, 12.5)
    6 = @(32, 2, 94.12,             5.88,  7.6, 0, 0)
    7 = @(27, 9, 75,                25,    7.2, 3, 8.33)
}

$sheetNames = @("1er Parcial", "3er Parcial")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $vals = $updates[$row]

        $ws.Cells.Item($row, 5).Value  = $vals[0]   # E - Aprobados
        $ws.Cells.Item($row, 6).Value  = $vals[1]   # F - Reprobados
        $ws.Cells.Item($row, 7).Value  = $vals[2]   # G - Por_Apro
        $ws.Cells.Item($row, 8).Value  = $vals[3]   # H - Por_Repro
        $ws.Cells.Item($row, 9).Value  = $vals[4]   # I - Promedio
        $ws.Cells.Item($row, 10).Value = $vals[5]   # J - Blancos
        $ws.Cells.Item($row, 11).Value = $vals[6]   # K - Por_Blan
    }
}
